$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "1/6/2020"
$ws.Range("B15").Value = "Basic java concepts"
$ws.Range("D15").Value = "Example programs"

$ws.Range("G11").Select()
